$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-170 down to 76-171.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with its data.
$ws.Cells.Item(75, 1).Value = 11
$ws.Cells.Item(75, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(75, 3).Value = "Bíobío"
$ws.Cells.Item(75, 4).Value = 45118
$ws.Cells.Item(75, 5).Value = 8
$ws.Cells.Item(75, 6).Value = 100112001
$ws.Cells.Item(75, 7).Value = "Berenjena"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 50
$ws.Cells.Item(75, 11).Value = 9000
$ws.Cells.Item(75, 12).Value = 9000
$ws.Cells.Item(75, 13).Value = 9000
$ws.Cells.Item(75, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(75, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(75, 16).Value = 180
$ws.Cells.Item(75, 17).Value = 50
$ws.Cells.Item(75, 18).Value = "Hortaliza"
